$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.446.53'

$ws.Range("E2").Value = '  +9.17%  '

$ws.Range("D3").Value = '1.679.12'

$ws.Range("E3").Value = '  +4.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("B5").Value = 'BNB'

$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.99'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +1.39%  '

$ws.Range("B6").Value = 'USDC'

$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3448'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  +1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.22'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +13.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.182'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  +3.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07273'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9998'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.44'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.141'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +3.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.751'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  +1.92%  '

$ws.Range("B16").Value = 'ShibaInu'

$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001111'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  +2.48%  '

$ws.Range("B17").Value = 'WrappedEther'

$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'

$ws.Range("D17").Value = '1.664.21'

$ws.Range("E17").Value = '  +4.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9977'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06718'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '81.35'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +4.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.48'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +2.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.113'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +1.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  +1.71%  '

$ws.Range("D24").Value = '24.400.84'

$ws.Range("E24").Value = '  +8.74%  '

$ws.Range("E25").Value = '  +1.82%  '

$ws.Range("B26").Value = 'LidoDAOToken'

$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.679'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  +6.08%  '

$ws.Range("B27").Value = 'LEO'

$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.359'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -11.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.44'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +1.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.56'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").Value = '1.861.58'

$ws.Range("E30").Value = '  +4.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.90'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  +4.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.375'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  +4.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.031'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  -3.64%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9755'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +2.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08474'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  +2.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.686'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +3.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '12.48'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  +5.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06497'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +6.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.960'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +4.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.362'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  +1.72%  '

$ws.Range("E41").Value = '  +5.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.271'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.2118'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  +4.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6194'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +4.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9977'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +0.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.774'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -1.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.07'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5957'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +4.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.05'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  +3.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07217'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +5.86%  '
